# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.581.96"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.999.31"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.014"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +1.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "329.97"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.85%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.012"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4993"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.17%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4222"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.12%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "53.73"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.09000"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.42%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.117"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.41%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "23.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.51%  "
$ws.Range("D13").Value = "2.021.73"
$ws.Range("E13").Value = "  -0.70%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "8.048"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -7.21%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.466"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -6.19%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.014"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "93.91"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -7.00%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001111"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.74%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06665"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "19.71"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.955"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -6.29%  "
$ws.Range("D23").Value = "29.597.04"
$ws.Range("E23").Value = "  -2.23%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.96"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.30%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.292"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "159.32"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.67"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.89%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.438"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.39%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.297"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -8.61%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "128.23"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.050"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -7.03%  "
$ws.Range("E32").Value = "  -5.06%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.573"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.58%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.828"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.13%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.803"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.35%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02465"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.08%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "9.349"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -7.42%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.305"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.82%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.06344"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -6.26%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6556"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.14%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "11.67"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.2049"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.38%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.011"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.88%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6328"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.97%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.49"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.17%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.194"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -6.07%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.304"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.80%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.514"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.09%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00000000336"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.64%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06984"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.18%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.123"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.25%  "
